# Apply scenario update to NCAP_BND sheet: scale a block of E-column
# figures down to 3/5 of their original values (entered as live formulas
# so Excel recomputes them), and set the single plain-value change in E32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NCAP_BND")

# Rows 12-21: replace literal numbers with "=<orig>*3/5" formulas.
$ws.Range("E12").Formula = "=13*3/5"
$ws.Range("E13").Formula = "=15*3/5"
$ws.Range("E14").Formula = "=17*3/5"
$ws.Range("E15").Formula = "=19*3/5"
$ws.Range("E16").Formula = "=22*3/5"
$ws.Range("E17").Formula = "=8.4*3/5"
$ws.Range("E18").Formula = "=12*3/5"
$ws.Range("E19").Formula = "=18*3/5"
$ws.Range("E20").Formula = "=25*3/5"
$ws.Range("E21").Formula = "=33*3/5"

# Row 32 is a direct value edit (40 -> 30), not a formula.
$ws.Range("E32").Value = 30

# Rows 33-36: same *3/5 scaling treatment as above.
$ws.Range("E33").Formula = "=90*3/5"
$ws.Range("E34").Formula = "=130*3/5"
$ws.Range("E35").Formula = "=180*3/5"
$ws.Range("E36").Formula = "=230*3/5"

# Update the window selection to match the saved state (cursor moved to E11,
# top-left cell reset to default A1 view).
$ws.Range("E11").Select()
